$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 31: Jane Smith
$ws.Cells.Item(31, 1).Value = 110030
$ws.Cells.Item(31, 2).Value = 9317596768
$ws.Cells.Item(31, 3).Value = "Jane Smith"
$ws.Cells.Item(31, 4).Value = "jane.smith@xyz.com"
$ws.Cells.Item(31, 5).Value = 818876432
$ws.Cells.Item(31, 6).Value = "ACT"
$ws.Cells.Item(31, 7).Value = "eng"
$ws.Cells.Item(31, 8).Value = "PWD"
$ws.Cells.Item(31, 9).Value = $true
$ws.Cells.Item(31, 10).Value = "superadmin"
$ws.Cells.Item(31, 11).Value = "now()"

# Row 32: John Doe
$ws.Cells.Item(32, 1).Value = 110031
$ws.Cells.Item(32, 2).Value = 9317596767
$ws.Cells.Item(32, 3).Value = "John Doe"
$ws.Cells.Item(32, 4).Value = "john.doe@xyz.com"
$ws.Cells.Item(32, 5).Value = 818876431
$ws.Cells.Item(32, 6).Value = "ACT"
$ws.Cells.Item(32, 7).Value = "eng"
$ws.Cells.Item(32, 8).Value = "PWD"
$ws.Cells.Item(32, 9).Value = $true
$ws.Cells.Item(32, 10).Value = "superadmin"
$ws.Cells.Item(32, 11).Value = "now()"

$ws.Range("E28").Select()
